# Update "想去人数" (F column) values per the source-data refresh (commit 456a3b4)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7266
$ws1.Range("F5").Value = 23
$ws1.Range("F6").Value = 568
$ws1.Range("F7").Value = 192
$ws1.Range("F11").Value = 61
$ws1.Range("F12").Value = 221
$ws1.Range("F13").Value = 17
$ws1.Range("F14").Value = 467
$ws1.Range("F16").Value = 1870
$ws1.Range("F17").Value = 50
$ws1.Range("F18").Value = 45
$ws1.Range("F19").Value = 3802
$ws1.Range("F21").Value = 253
$ws1.Range("F23").Value = 44
$ws1.Range("F24").Value = 2
$ws1.Range("F25").Value = 37
$ws1.Range("F26").Value = 2478
$ws1.Range("F27").Value = 24
$ws1.Range("F28").Value = 314
$ws1.Range("F30").Value = 7
$ws1.Range("F31").Value = 43
$ws1.Range("F32").Value = 9
$ws1.Range("F38").Value = 21
$ws1.Range("F39").Value = 1478
$ws1.Range("F40").Value = 161

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7266
$ws4.Range("F5").Value = 23
$ws4.Range("F7").Value = 568
$ws4.Range("F8").Value = 192
$ws4.Range("F12").Value = 61
$ws4.Range("F13").Value = 221
$ws4.Range("F14").Value = 17
$ws4.Range("F15").Value = 467
$ws4.Range("F17").Value = 1870
$ws4.Range("F18").Value = 50
$ws4.Range("F19").Value = 45
$ws4.Range("F20").Value = 3802
$ws4.Range("F22").Value = 253
$ws4.Range("F24").Value = 44
$ws4.Range("F25").Value = 2
$ws4.Range("F26").Value = 37
$ws4.Range("F27").Value = 2479
$ws4.Range("F28").Value = 24
$ws4.Range("F29").Value = 314
$ws4.Range("F31").Value = 7
$ws4.Range("F32").Value = 43
$ws4.Range("F33").Value = 9
$ws4.Range("F39").Value = 21
$ws4.Range("F40").Value = 1478
$ws4.Range("F41").Value = 161

